$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the diff.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the workbook's original inlineStr/text cells) instead of
# auto-converting numeric-looking strings (e.g. '231.57') into numbers;
# resetting the Style to 'Normal' afterwards avoids leaving a stray
# quote-prefix style on the cell so formatting matches the source file.

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.Value = "'37.127.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.Value = "'  +1.50%  "
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.Value = "'2.050.11"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.Value = "'  +0.08%  "
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4, 5)
$c.Value = "'  +0.00%  "
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.Value = "'231.57"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.Value = "'  +0.00%  "
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.622"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.Value = "'  +3.49%  "
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7, 5)
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.Value = "'57.15"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.Value = "'  +3.27%  "
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.381"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.Value = "'  +2.62%  "
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.Value = "'57.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.Value = "'  +0.23%  "
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.0755"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.Value = "'  +0.90%  "
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.102"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.Value = "'  +1.22%  "
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.Value = "'2.351.46"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.Value = "'  +0.54%  "
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.Value = "'14.26"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.Value = "'  -0.97%  "
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.Value = "'20.77"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.Value = "'  +3.19%  "
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.Value = "'0.770"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.Value = "'  +1.03%  "
$c.Style = "Normal"

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.Value = "'5.13"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.Value = "'  -0.03%  "
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.Value = "'2.046.42"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.Value = "'  +0.21%  "
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.Value = "'37.094.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.Value = "'  +0.83%  "
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.Value = "'6.27"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.Value = "'  +9.08%  "
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.Value = "'68.77"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.Value = "'  +1.62%  "
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.Value = "'0.0₃0807"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.Value = "'  +1.32%  "
$c.Style = "Normal"

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.Value = "'224.74"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.Value = "'  +1.72%  "
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 5)
$c.Value = "'  +0.09%  "
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.43"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.Value = "'  +1.28%  "
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.Value = "'2.37"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.Value = "'166.09"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.Value = "'  +1.84%  "
$c.Style = "Normal"

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.Value = "'1.45"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.Value = "'  +7.24%  "
$c.Style = "Normal"

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.Value = "'8.75"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.Value = "'  +0.37%  "
$c.Style = "Normal"

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.Value = "'18.98"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.Value = "'  -0.19%  "
$c.Style = "Normal"

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.Value = "'0.125"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.Value = "'  -3.14%  "
$c.Style = "Normal"

# Row 32
$c = $ws.Cells.Item(32, 5)
$c.Value = "'  -0.25%  "
$c.Style = "Normal"

# Row 33
$c = $ws.Cells.Item(33, 5)
$c.Value = "'  +1.05%  "
$c.Style = "Normal"

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.0613"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.Value = "'  +1.61%  "
$c.Style = "Normal"

# Row 35
$c = $ws.Cells.Item(35, 2)
$c.Value = "'InternetComputer(DFINITY)"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 3)
$c.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.Value = "'4.56"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.Value = "'  +6.63%  "
$c.Style = "Normal"

# Row 36
$c = $ws.Cells.Item(36, 2)
$c.Value = "'LidoDAOToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 3)
$c.Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.Value = "'2.52"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.Value = "'  -0.77%  "
$c.Style = "Normal"

# Row 37
$c = $ws.Cells.Item(37, 5)
$c.Value = "'  -0.08%  "
$c.Style = "Normal"

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.Value = "'1.74"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.Value = "'  -0.72%  "
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39, 5)
$c.Value = "'  -0.51%  "
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.Value = "'5.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.Value = "'  -2.31%  "
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41, 5)
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42, 5)
$c.Value = "'  +1.81%  "
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.474.93"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.Value = "'96.42"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.Value = "'  +3.27%  "
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45, 2)
$c.Value = "'TrustWalletToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 3)
$c.Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'1.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.Value = "'  +4.09%  "
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46, 2)
$c.Value = "'Cronos"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.0924"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.Value = "'  -1.41%  "
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.Value = "'0.0209"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.Value = "'  +2.65%  "
$c.Style = "Normal"

# Row 48
$c = $ws.Cells.Item(48, 2)
$c.Value = "'ARBITRUM"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'1.02"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.Value = "'  +0.68%  "
$c.Style = "Normal"

# Row 49
$c = $ws.Cells.Item(49, 2)
$c.Value = "'InjectiveProtocol"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 3)
$c.Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.Value = "'15.11"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.Value = "'  -3.54%  "
$c.Style = "Normal"

# Row 50
$c = $ws.Cells.Item(50, 2)
$c.Value = "'FraxShare"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 3)
$c.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.Value = "'7.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.Value = "'  +3.23%  "
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 2)
$c.Value = "'MXToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 3)
$c.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.Value = "'2.93"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.Value = "'  +1.22%  "
$c.Style = "Normal"
